$d = $word.ActiveDocument

function FindReplace($old, $new) {
    $ok = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $ok) {
        Write-Output "MISSING: $old"
    }
    return $ok
}

function SplitRange($rng) {
    # Toggling a character-level property and back forces the range to
    # become its own run without changing its visible formatting.
    $rng.Font.Bold = 1
    $rng.Font.Bold = 0
}

# =========================================================================
# 1. Bulleted "Profile" list: numId 4 -> 5 (shared new list) for the three
#    bullet paragraphs (Over 20 years / Strong expertise / Solid RDBMS).
# =========================================================================

$p5 = $d.Paragraphs.Item(5)
$p7 = $d.Paragraphs.Item(7)
if (($p5.Range.Text -like "*Over 20 years*") -and ($p7.Range.Text -like "*Solid RDBMS*")) {
    $listRange = $d.Range($p5.Range.Start, $p7.Range.End)
    $gallery = $d.ListGalleries.Item(1)
    $template = $gallery.ListTemplates.Item(1)
    $listRange.ListFormat.ApplyListTemplate($template)
    Write-Output "numbering updated"
} else {
    Write-Output "profile paragraphs not where expected"
}

# =========================================================================
# 2. " Spring Core 2/3/4, MVC, ..." -> add "Spring Boot 1.5, " prefix
# =========================================================================

FindReplace " Spring Core 2/3/4, MVC, JDBC Template, Security, Web Services, Rest Template, AMQP, Spring Data, Spring Data " " Spring Boot 1.5, Spring Core 2/3/4, MVC, JDBC Template, Security, Web Services, Rest Template, AMQP, Spring Data, Spring Data "

# =========================================================================
# 3. "Agile, SCRUM, Test Driven Development " -> insert "Kanban, " as its
#    own run (mirrors the spell-checked insertion in the original edit).
# =========================================================================

$r = $d.Content
if ($r.Find.Execute("Agile, SCRUM, ")) {
    $r.Collapse(0)
    $insStart = $r.Start
    $r.InsertAfter("Kanban, ")
    $kanbanRange = $d.Range($insStart, $insStart + 6)
    SplitRange $kanbanRange
    Write-Output "Kanban inserted"
} else {
    Write-Output "MISSING: Agile, SCRUM, "
}

# =========================================================================
# 4. "Tomcat 7/8, Oracle 11g, Azure and private VMware cloud platforms, ..."
# =========================================================================

$ok = (FindReplace "Tomcat 7/8, Oracle 11g, Azure and private VMware cloud platforms, Linux, Windows 7" "Tomcat 7/8, Oracle 11g, Azure, Pivotal Cloud Foundry and private VMware cloud platforms, Linux, Windows 7")
if ($ok) {
    $r2 = $d.Content
    if ($r2.Find.Execute("cloud platforms, Linux, Windows 7")) {
        $r2.Collapse(1)
        $cloudStart = $r2.Start
        $cloudRange = $d.Range($cloudStart, $cloudStart + 5)
        $cloudRange.Text = "OneCloud"
        $oneCloudRange = $d.Range($cloudStart, $cloudStart + 8)
        SplitRange $oneCloudRange
        Write-Output "Tomcat OneCloud split"
    } else {
        Write-Output "MISSING: cloud platforms, Linux, Windows 7"
    }
}

# =========================================================================
# 5. ", Bamboo, Oracle " -> ", Bamboo, Concourse, Oracle "
# =========================================================================

FindReplace ", Bamboo, Oracle " ", Bamboo, Concourse, Oracle "

# =========================================================================
# 6. Big McKesson bullet paragraph rewrite.
# =========================================================================

$old1 = "Created two Spring Boot based micro services providing new patient and formulary functionality, enabling independent scalability and cloud deployment. "
$new1 = "Created several Spring Boot based micro services providing new patient, formulary, order validation and submission functionality. Used private Vmware OneCloud, Azure and Pivotal Cloud Foundry to enable application independent scalability and cloud deployment. "
$ok1 = (FindReplace $old1 $new1)

if ($ok1) {
    # split "order"
    $r3 = $d.Content
    if ($r3.Find.Execute("order validation")) {
        $orderStart = $r3.Start
        $orderRange = $d.Range($orderStart, $orderStart + 5)
        SplitRange $orderRange
    }

    # split "Vmware"
    $r4 = $d.Content
    if ($r4.Find.Execute("Used private Vmware OneCloud")) {
        $vStart = $r4.Start + ("Used private ").Length
        $vRange = $d.Range($vStart, $vStart + 6)
        SplitRange $vRange

        # split "OneCloud" right after Vmware + " "
        $oStart = $vStart + 7
        $oRange = $d.Range($oStart, $oStart + 8)
        SplitRange $oRange
    }

    # split " Completed " into " " + "Completed "
    $r5 = $d.Content
    if ($r5.Find.Execute(" Completed ")) {
        $cStart = $r5.Start + 1
        $cRange = $d.Range($cStart, $cStart + ("Completed ").Length)
        SplitRange $cRange
    }

    # tail replace + split final "OneCloud"
    $old2 = " and Cloud Foundry POCs using Azure and private VMware cloud platforms helping with the cloud deployment strategy selection."
    $new2 = " and Pivotal Cloud Foundry POCs using Azure and private VMware OneCloud platforms helping with the cloud deployment strategy selection."
    if ((FindReplace $old2 $new2)) {
        $r6 = $d.Content
        if ($r6.Find.Execute("private VMware OneCloud platforms")) {
            $fStart = $r6.Start + ("private VMware ").Length
            $fRange = $d.Range($fStart, $fStart + 8)
            SplitRange $fRange
        }
    }
    Write-Output "McKesson paragraph rewritten"
}

# =========================================================================
# 7. Education section + footer date
# =========================================================================

FindReplace "Novosibirsk Technical University, Novosibirsk, Russia" "Novosibirsk State Technical University, Novosibirsk, Russia"

FindReplace "Last modified: 2017/06/03" "Last modified: 2017/12/26"

Write-Output "all edits applied"
